$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new rows describing the new "dash_configs" entries
# (multiple alias support for the dashboard config, per commit #200)
$ws.Cells.Item(79, 1).Value2 = 10.1
$ws.Cells.Item(79, 2).Value = "dash_configs"
$ws.Cells.Item(79, 3).Value = "dash_configs"
$ws.Cells.Item(79, 4).Value = "alias_dash_lista"

$ws.Cells.Item(80, 1).Value2 = 10.2
$ws.Cells.Item(80, 2).Value = "dash_configs"
$ws.Cells.Item(80, 3).Value = "dash_configs"
$ws.Cells.Item(80, 4).Value = "alias_data_lista"

$ws.Cells.Item(81, 1).Value2 = 10.3
$ws.Cells.Item(81, 2).Value = "dash_configs"
$ws.Cells.Item(81, 3).Value = "dash_configs"
$ws.Cells.Item(81, 4).Value = "alias_insumos_lista"

# Re-apply the sort on column A (the whole table is kept sorted by "orden")
$sort = $ws.Sort
$sort.SortFields.Clear()
$key1 = $ws.Range("A2:A78")
$sort.SortFields.Add($key1, $null, 1, $null)
$sort.SetRange($ws.Range("A2:J78"))
$sort.Header = 2
$sort.Apply()

# Update the sheet view: select A6:XFD8 and scroll back to the top
$ws.Activate()
[void]$ws.Range("A6:XFD8").Select()
